$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be stored as text so numeric-looking price strings
# (e.g. "1.00", "8.31") are preserved verbatim instead of becoming numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.612.60"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "3.783.41"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "597.57"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").Value = "164.58"
$ws.Range("E6").Value = "  -1.58%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "0.515"
$ws.Range("E8").Value = "  -0.84%  "
$ws.Range("E9").Value = "  -1.01%  "
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("D11").Value = "6.42"
$ws.Range("E11").Value = "  +1.53%  "
$ws.Range("D12").Value = "0.0000248"
$ws.Range("E12").Value = "  -2.19%  "
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("D14").Value = "4.421.90"
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").Value = "3.794.29"
$ws.Range("E15").Value = "  -1.48%  "
$ws.Range("D16").Value = "67.685.65"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "18.26"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("E18").Value = "  +1.72%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "461.76"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("E21").Value = "  -2.39%  "
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").Value = "82.61"
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("D24").Value = "0.0000145"
$ws.Range("E24").Value = "  -6.77%  "
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").Value = "2.09"
$ws.Range("E26").Value = "  -1.37%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").Value = "3.932.41"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("D30").Value = "7.40"
$ws.Range("E30").Value = "  +2.18%  "
$ws.Range("E31").Value = "  -4.43%  "
$ws.Range("E32").Value = "  -2.11%  "
$ws.Range("E33").Value = "  -2.26%  "
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").Value = "8.97"
$ws.Range("E35").Value = "  -1.27%  "
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").Value = "0.985"
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("E40").Value = "  -5.00%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D43").Value = "43.91"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").Value = "47.50"
$ws.Range("E44").Value = "  -1.08%  "
$ws.Range("D45").Value = "0.297"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("D46").Value = "151.67"
$ws.Range("E46").Value = "  +1.49%  "
$ws.Range("D47").Value = "8.31"
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "1.36"
$ws.Range("E48").Value = "  +7.03%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "27.11"
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").Value = "394.05"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "1.85"
$ws.Range("E51").Value = "  +1.50%  "
